# 8.10.1.xlsx update: add a new "2022" column (S) of data and revise the
# already-published 2021 (R) figures that depended on the (now corrected)
# adult-population total for that year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring column S (the new year's column) up to the same formatting as the
# existing last-year column (R) for the header + data rows, before filling
# in values.
$ws.Range("R3:R8").Copy($ws.Range("S3:S8"))

# --- Row 3: year headers ---
$ws.Range("S3").Value2 = 2022

# --- Row 4: a) branches per 100 000 adults ---
# R4 used to be the live formula =R6/R8*100000; the 2021 population total
# (R8) was revised, so the figure is now stored as a plain recalculated
# value instead of a formula.
$ws.Range("R4").Value2 = 6.9132648934880807
$ws.Range("S4").Value2 = 6.9031689452913012

# --- Row 5: b) ATMs per 100 000 adults ---
# Same story as row 4 (was =R7/R8*100000).
$ws.Range("R5").Value2 = 42.321589572314856
$ws.Range("S5").Value2 = 44.306188104841333

# --- Row 6: total branches of commercial banks ---
$ws.Range("S6").Value2 = 318

# --- Row 7: total ATMs ---
$ws.Range("S7").Value2 = 2041

# --- Row 8: adult resident population (revised 2021 value + new 2022 value) ---
$ws.Range("R8").Value2 = 4513063
$ws.Range("S8").Value2 = 4606580

# Move the active selection, matching the author's final cursor position.
$ws.Range("R13").Select() | Out-Null
